# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" quarter snapshot ahead of the existing
# "2022-Q2" / "2021-Q3" sheets, pushing the older quarters down, and
# records the new quarter's numbers in the "总计" (totals) roll-up sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" roll-up sheet: insert a row for "2022-Q2" (it is no longer
#    the latest quarter) ahead of "2021-Q3", and relabel the former
#    "2022-Q2" row as "2022-Q3" (the new latest quarter keeps the same
#    aggregate counts: 1 fund holder, 0.02 (亿元) held).
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room: push row 3 ("2021-Q3") down to row 4.
$total.Rows.Item(3).Insert()

# Copy row 2's formatting down into the freshly inserted row 3 so the
# index cell keeps the same border/bold/center style as its neighbours.
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)

# Row 2 becomes the new latest quarter, "2022-Q3" (counts unchanged).
$total.Range("B2").Value = "2022-Q3"

# Row 3: the previous quarter, "2022-Q2", with the counts it had when
# it was still the latest quarter.
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.02

# Row 4: "2021-Q3" shifted down one row; only its running index changes.
$total.Range("A4").Value = 2

# ------------------------------------------------------------------
# 2) Duplicate the current "2022-Q2" sheet (with its original holdings
#    data) and place the copy right before "2021-Q3" -- this becomes
#    the preserved "2022-Q2" sheet once the original is promoted to
#    "2022-Q3" below. The original is renamed out of the way first so
#    the copy can keep the plain "2022-Q2" name without colliding.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3last = $wb.Worksheets.Item("2021-Q3")
$q2.Copy($q3last)
$q2.Name = "2022-Q3"
$wb.Worksheets.Item("2022-Q2 (2)").Name = "2022-Q2"

# ------------------------------------------------------------------
# 3) Turn the original sheet (now "2022-Q3") into the new quarter by
#    updating its fund-holding figures in place.
#    Text-formatted numeric-looking values are written with a leading
#    apostrophe (forces text storage) and then restyled back to the
#    workbook's default "Normal" style so no stray number format is
#    left behind.
# ------------------------------------------------------------------
$q2.Range("D2").Value = "'0.51"
$q2.Range("D2").Style = "Normal"

$q2.Range("E2").Value = "'79.37"
$q2.Range("E2").Style = "Normal"

$q2.Range("F2").Value = "'4.06"
$q2.Range("F2").Style = "Normal"

$q2.Range("G2").Value = "'0.0207"
$q2.Range("G2").Style = "Normal"

$q2.Range("H2").Value = 8

# Keep "2021-Q3" as the selected/active tab, matching its state before
# this edit (the newly duplicated "2022-Q2" sheet should not steal the
# active-tab flag). Re-fetched by name since the $q3last handle can go
# stale once other sheets have been inserted/renamed around it.
$wb.Worksheets.Item("2021-Q3").Activate()
